# Update gh-pages to output generated at 456a3b4
# Apply updated "F" column (售票量/sales count) values across the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

function Apply-FColumnUpdates($SheetName, $RowValues) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $RowValues[$row]
    }
}

# ----- 展览 (sheet1) -----
$exhibition = @{
    3  = 311
    4  = 2964
    7  = 2314
    8  = 1656
    10 = 849
    11 = 117
    12 = 8
    13 = 2655
    15 = 1506
    16 = 7029
    18 = 7183
    20 = 5399
    21 = 3102
    22 = 3470
    23 = 227
    25 = 1863
    26 = 80
    28 = 876
    29 = 21
    30 = 176
    31 = 38
    32 = 2403
    33 = 1160
    34 = 2664
    35 = 20
    37 = 167
    38 = 383
    39 = 1057
    41 = 473
    42 = 520
}
Apply-FColumnUpdates "展览" $exhibition

# ----- 演出 (sheet2) -----
$performance = @{
    6 = 1
    7 = 35
    8 = 213
}
Apply-FColumnUpdates "演出" $performance

# ----- 全部类型 (sheet4) -----
$allTypes = @{
    4  = 311
    6  = 2964
    8  = 2314
    9  = 1656
    11 = 849
    12 = 117
    13 = 35
    14 = 2655
    15 = 1506
    16 = 213
    19 = 7029
    21 = 7183
    23 = 5399
    24 = 3102
    25 = 3470
    27 = 227
    29 = 1863
    33 = 876
    34 = 21
    35 = 176
    36 = 38
    37 = 2403
    38 = 1160
    40 = 2664
    41 = 20
    43 = 167
    45 = 383
    46 = 1057
    48 = 473
    49 = 520
}
Apply-FColumnUpdates "全部类型" $allTypes
